$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.920.58'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.874.07'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'0.7431"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.93%  '
$ws.Range('D6').Value = "'242.67"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.3152"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').Value = "'0.07220"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = "'24.65"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.94%  '
$ws.Range('D11').Value = "'0.08330"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('D12').Value = "'0.7503"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.379"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.867.94'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = "'92.25"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.45%  '
$ws.Range('D16').Value = "'6.120"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '29.920.87'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').Value = "'247.08"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('E19').Value = '  -1.41%  '
$ws.Range('D20').Value = "'0.000007841"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.144.93'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.0000"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').Value = "'8.008"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').Value = "'1.000"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = "'0.1543"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.97%  '
$ws.Range('D26').Value = "'9.293"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('D27').Value = "'165.48"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('D29').Value = "'2.019"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('D30').Value = "'1.501"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('D31').Value = "'4.592"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('D32').Value = "'1.537"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = "'4.219"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.02%  '
$ws.Range('D34').Value = "'0.05337"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('D35').Value = "'1.234"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('D36').Value = "'0.7499"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').Value = "'2.702"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('D39').Value = "'0.01965"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').Value = "'2.754"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('D41').Value = "'0.4518"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').Value = '1.115.59'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('D43').Value = "'6.132"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('D44').Value = "'72.44"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').Value = "'0.8629"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').Value = "'104.35"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = "'1.864"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('D49').Value = "'7.623"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').Value = "'9.515"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('D51').Value = '2.037.30'
$ws.Range('E51').Value = '  +0.36%  '
